# B6-PowerPoint.pptx edit
#
# The authoritative change captured by this commit is a table-style
# re-skin: every table in the deck that used the custom "Table_0" style
# ({D16B7D89-1DE2-4527-B30B-859C5D4928BF}) is switched to the built-in
# PowerPoint table style {6728B564-69C5-48F8-AC21-8E34BD33B002}.
#
# Walk every slide/shape looking for tables and re-apply the new style
# wherever the old style id is found (this covers slides 14, 15 and 16
# without hard-coding shape indices).

$OLD_STYLE = "{D16B7D89-1DE2-4527-B30B-859C5D4928BF}"
$NEW_STYLE = "{6728B564-69C5-48F8-AC21-8E34BD33B002}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $OLD_STYLE) {
                $tbl.ApplyStyle($NEW_STYLE)
            }
        }
    }
}
